# Chapter 6 / Example 5 - add "最大利润" (Max profit) and "最小利润" (Min profit)
# columns (I, J) to each of the five product worksheets, using Excel's
# MAX()/MIN() values over the per-sheet "销售利润" (profit) column (H).
#
# Currency number format used throughout the workbook (numFmtId 7 equivalent).
$currencyFmt = '"¥"#,##0.00;"¥"\-#,##0.00'

$wb = $excel.ActiveWorkbook

$maxHeader = "最大利润"
$minHeader = "最小利润"

# ---------------------------------------------------------------------
# Step 1: write the "最大利润" header (column I) on every sheet. Each
# sheet's first row already carries its own style (row-level or via the
# column default), so a plain value write naturally reuses/matches the
# existing header formatting.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 5; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("I1").Value = $maxHeader
}

# ---------------------------------------------------------------------
# Step 2: write the "最小利润" header (column J). On sheets 1-4 this was
# typed/bolded and given the currency format (creating a new bold +
# centered + currency style); on sheet 5 it only picked up the existing
# currency/centered style already used for the numeric columns.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("J1").Value = $minHeader
    $ws.Range("J1").Font.Bold = $true
    $ws.Range("J1").HorizontalAlignment = -4108
    $ws.Range("J1").NumberFormat = $currencyFmt
}

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("J1").Value = $minHeader
$ws5.Range("J1").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------
# Step 3: fill in the MAX/MIN profit figures (row 2) for each sheet,
# reading straight off the existing H column ("销售利润") so the numbers
# always reflect that sheet's own data.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 5; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $maxProfit = $ws.Application.WorksheetFunction.Max($ws.Range("H2:H13"))
    $minProfit = $ws.Application.WorksheetFunction.Min($ws.Range("H2:H13"))

    $ws.Range("I2").Value = $maxProfit
    $ws.Range("J2").Value = $minProfit
    $ws.Range("I2").NumberFormat = $currencyFmt
    $ws.Range("J2").NumberFormat = $currencyFmt
}

Write-Output "added max/min profit columns to all sheets"
